# Updated cryptos list on Mon Oct  2 13:23:46 UTC 2023 with GitHub Actions
#
# Applies the latest price / 1h-volume refresh to the cryptos worksheet.
# Some "Price" (column D) values look like plain numbers (e.g. "24.05"),
# which Excel would normally auto-convert to a Number when assigned to a
# General-formatted cell. The source cells are text, so those particular
# assignments are entered with a leading apostrophe (quote-prefix) to force
# text entry, then the cell style is reset back to "Normal" so no stray
# number-format / quote-prefix styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.Value = "'" + $Text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.361.91"
$ws.Range("E2").Value = "  +4.20%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.728.71"
$ws.Range("E3").Value = "  +2.43%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.19%  "

# Row 5 - BNB
Set-TextValue "D5" "219.11"
$ws.Range("E5").Value = "  +1.42%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.57%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.22%  "

# Row 8 - Solana
Set-TextValue "D8" "24.05"
$ws.Range("E8").Value = "  +2.72%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.268"
$ws.Range("E9").Value = "  +2.17%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0636"
$ws.Range("E10").Value = "  +1.32%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0895"
$ws.Range("E11").Value = "  +0.60%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.972.95"

# Row 13 - WrappedEther
Set-TextValue "D13" "1.727.42"
$ws.Range("E13").Value = "  +2.39%  "

# Row 14 - Polkadot
Set-TextValue "D14" "4.24"
$ws.Range("E14").Value = "  +0.86%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.566"
$ws.Range("E15").Value = "  +2.26%  "

# Row 16 - Litecoin
Set-TextValue "D16" "67.63"
$ws.Range("E16").Value = "  +0.33%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "28.326.03"
$ws.Range("E17").Value = "  +4.06%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "246.37"
$ws.Range("E18").Value = "  +4.14%  "

# Row 19 - ShibaInu
Set-TextValue "D19" "0.0₃0753"
$ws.Range("E19").Value = "  +1.13%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -1.95%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.23%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.61"
$ws.Range("E22").Value = "  +1.25%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.66"
$ws.Range("E23").Value = "  +0.22%  "

# Row 24 - Toncoin
Set-TextValue "D24" "2.06"
$ws.Range("E24").Value = "  -2.10%  "

# Row 25 - Monero
Set-TextValue "D25" "149.37"
$ws.Range("E25").Value = "  +1.39%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.45"
$ws.Range("E26").Value = "  +1.88%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "16.61"
$ws.Range("E27").Value = "  +1.11%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +0.37%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  -0.29%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +2.75%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.72%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.42"
$ws.Range("E32").Value = "  +0.72%  "

# Row 33 - Maker
Set-TextValue "D33" "1.484.73"
$ws.Range("E33").Value = "  -4.30%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +0.43%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "1.64"
$ws.Range("E35").Value = "  -1.75%  "

# Row 36 - ARBITRUM
Set-TextValue "D36" "0.978"
$ws.Range("E36").Value = "  +3.26%  "

# Row 37 - was ImmutableX, now HuobiToken (rows 37/38 swapped order)
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D37" "2.41"
$ws.Range("E37").Value = "  +0.35%  "

# Row 38 - was HuobiToken, now ImmutableX
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D38" "0.602"
$ws.Range("E38").Value = "  -0.21%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +1.29%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  +0.24%  "

# Row 41 - Aave
Set-TextValue "D41" "69.73"
$ws.Range("E41").Value = "  +0.76%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.27%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -1.96%  "

# Row 44 - RocketPoolETH
Set-TextValue "D44" "1.877.86"
$ws.Range("E44").Value = "  +2.21%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  +1.04%  "

# Row 46 - TrustWalletToken
Set-TextValue "D46" "0.806"
$ws.Range("E46").Value = "  +2.28%  "

# Row 47 - RenderToken
Set-TextValue "D47" "1.73"
$ws.Range("E47").Value = "  +7.67%  "

# Row 48 - was Quant, now BabyDogeCoin (rows 48/49 swapped order)
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.0₆0114"
$ws.Range("E48").Value = "  +3.89%  "

# Row 49 - was BabyDogeCoin, now Quant
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D49" "90.30"
$ws.Range("E49").Value = "  -1.02%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "8.14"
$ws.Range("E50").Value = "  -3.18%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -0.68%  "
